$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Price (D) / Volume(1h) (E) figures, per-row, as published by the
# "Updated cryptos list" GitHub Action refresh.
$updates = @(
    @{ Row = 2; D = "27.024.97"; E = "  +0.67%  " },
    @{ Row = 3; D = "1.825.54"; E = "  +0.76%  " },
    @{ Row = 4; D = "1.006"; E = "  +0.46%  " },
    @{ Row = 5; D = "311.39"; E = "  +0.35%  " },
    @{ Row = 6; D = $null; E = "  +0.32%  " },
    @{ Row = 7; D = "0.4695"; E = "  -0.29%  " },
    @{ Row = 8; D = "0.3672"; E = "  -0.68%  " },
    @{ Row = 9; D = "0.07368"; E = "  +0.25%  " },
    @{ Row = 10; D = "0.8762"; E = "  +0.86%  " },
    @{ Row = 11; D = $null; E = "  -0.38%  " },
    @{ Row = 12; D = "1.842.88"; E = "  -1.60%  " },
    @{ Row = 13; D = "0.07314"; E = "  +3.44%  " },
    @{ Row = 14; D = "5.439"; E = "  +1.69%  " },
    @{ Row = 15; D = "6.528"; E = "  +0.38%  " },
    @{ Row = 16; D = "91.91"; E = "  +0.32%  " },
    @{ Row = 17; D = $null; E = "  +0.41%  " },
    @{ Row = 18; D = "0.000008749"; E = "  +0.60%  " },
    @{ Row = 19; D = $null; E = "  +0.24%  " },
    @{ Row = 20; D = "14.71"; E = "  +0.22%  " },
    @{ Row = 21; D = "27.033.91"; E = "  +0.55%  " },
    @{ Row = 22; D = "5.282"; E = "  -0.92%  " },
    @{ Row = 23; D = "10.64"; E = "  +0.95%  " },
    @{ Row = 24; D = "2.055.04"; E = "  -1.81%  " },
    @{ Row = 25; D = "1.894"; E = "  -0.11%  " },
    @{ Row = 26; D = "151.33"; E = "  -0.36%  " },
    @{ Row = 27; D = "18.46"; E = "  +0.32%  " },
    @{ Row = 28; D = "2.149"; E = "  +2.62%  " },
    @{ Row = 29; D = "5.241"; E = "  -0.87%  " },
    @{ Row = 30; D = "116.77"; E = "  +1.36%  " },
    @{ Row = 31; D = "0.08891"; E = "  -0.50%  " },
    @{ Row = 32; D = "0.7567"; E = "  +0.19%  " },
    @{ Row = 33; D = "1.163"; E = "  +1.26%  " },
    @{ Row = 34; D = "4.517"; E = "  +1.36%  " },
    @{ Row = 35; D = "2.932"; E = "  +0.10%  " },
    @{ Row = 36; D = $null; E = "  +0.34%  " },
    @{ Row = 37; D = "1.097"; E = "  +0.18%  " },
    @{ Row = 38; D = "0.05315"; E = "  +1.20%  " },
    @{ Row = 39; D = $null; E = "  +0.03%  " },
    @{ Row = 40; D = "2.979"; E = "  +2.47%  " },
    @{ Row = 41; D = "7.232"; E = "  +0.97%  " },
    @{ Row = 42; D = "2.385"; E = "  +1.64%  " },
    @{ Row = 43; D = "0.5311"; E = "  -0.28%  " },
    @{ Row = 44; D = "0.1658"; E = "  -0.05%  " },
    @{ Row = 45; D = "8.496"; E = "  +0.95%  " },
    @{ Row = 46; D = "0.4910"; E = "  -0.41%  " },
    @{ Row = 47; D = "10.50"; E = "  +1.55%  " },
    @{ Row = 48; D = $null; E = "  +0.31%  " },
    @{ Row = 49; D = "1.667"; E = "  -0.21%  " },
    @{ Row = 50; D = "103.38"; E = "  +0.35%  " },
    @{ Row = 51; D = $null; E = "  +0.35%  " }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($null -ne $u.D) {
        # Prefix with an apostrophe so Excel stores the price as literal
        # text (matches the source data, which includes values such as
        # "27.024.97" that are not valid numbers) instead of re-parsing it
        # as a number; then reset the style so no numeric/text format is
        # left behind on the cell.
        $ws.Range("D$row").Value = "'" + $u.D
        $ws.Range("D$row").Style = "Normal"
    }

    if ($null -ne $u.E) {
        $ws.Range("E$row").Value = $u.E
    }
}
